$d = $word.ActiveDocument

# The sentence below is rewritten to add " separados por espacio" in the
# middle of it. In the canonical OOXML this turns the paragraph's single
# <w:r> into three consecutive runs (all with identical run formatting):
#   1) "...riesgo ALTO"
#   2) " separados por espacio"
#   3) ", en caso de no haber ninguno devolver NA."
$beforeText = "Indicar el nombre de los municipios de los cuerpos de agua que tienen un nivel de riesgo ALTO"
$afterText  = ", en caso de no haber ninguno devolver NA."
$insertion  = " separados por espacio"
$fullOld    = $beforeText + $afterText

function Esc([string]$s) {
    return $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

$r = $d.Content
$found = $r.Find.Execute($fullOld, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Pull the live canonical XML for the story so we can read the exact
    # paragraph/run properties (paraId, rsids, numbering, color, ...) of the
    # sentence we just located, instead of inventing/guessing formatting.
    $full = $r.WordOpenXML
    $marker = $fullOld + "</w:t></w:r></w:p>"
    $markerIdx = $full.IndexOf($marker)

    if ($markerIdx -ge 0) {
        $prefix = $full.Substring(0, $markerIdx)

        # Paragraph opening tag (w14:paraId, rsids, ...) and its <w:pPr>.
        $pStart = $prefix.LastIndexOf("<w:p ")
        $pOpenEnd = $full.IndexOf(">", $pStart) + 1
        $pOpenTag = $full.Substring($pStart, $pOpenEnd - $pStart)

        $afterP = $full.Substring($pOpenEnd)
        if ($afterP.StartsWith("<w:pPr>")) {
            $pPrLen = $afterP.IndexOf("</w:pPr>") + "</w:pPr>".Length
            $pPr = $afterP.Substring(0, $pPrLen)
        } else {
            $pPr = ""
        }

        # Run formatting (<w:rPr>) of the run that currently holds the sentence.
        $runStart = $prefix.LastIndexOf("<w:r>")
        $afterR = $full.Substring($runStart)
        if ($afterR.StartsWith("<w:r><w:rPr>")) {
            $rPrLen = $afterR.IndexOf("</w:rPr>") + "</w:rPr>".Length - "<w:r>".Length
            $rPr = $afterR.Substring("<w:r>".Length, $rPrLen)
        } else {
            $rPr = ""
        }

        $run1 = "<w:r>$rPr<w:t>" + (Esc $beforeText) + "</w:t></w:r>"
        $run2 = "<w:r>$rPr<w:t xml:space=""preserve"">" + (Esc $insertion) + "</w:t></w:r>"
        $run3 = "<w:r>$rPr<w:t>" + (Esc $afterText) + "</w:t></w:r>"

        $newPara = "$pOpenTag$pPr$run1$run2$run3</w:p>"

        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
               "<w:body>$newPara</w:body>" +
               '</w:document></pkg:xmlData></pkg:part></pkg:package>'

        $r.InsertXML($xml)
        Write-Output "Paragraph updated (runs split in three)."
    } else {
        # Fallback: plain text replace (keeps a single run) if the live XML
        # could not be parsed the way we expected.
        $r2 = $d.Content
        $r2.Find.Execute($fullOld, $true, $false, $false, $false, $false, $true, 1, $false, $beforeText + $insertion + $afterText, 2)
        Write-Output "Paragraph updated (fallback plain replace)."
    }
} else {
    Write-Output "Target sentence not found; no changes made."
}
